# Updated cryptos list on Mon Mar  6 05:33:54 UTC 2023 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row, and
# swaps the VeChain / TrustWalletToken rows (37 and 38) including their
# name, link, price and volume values.
#
# Numeric-looking price strings are written with a leading apostrophe so
# Excel keeps them as literal text (matching the original inlineStr cells)
# instead of silently re-parsing them as numbers, which would otherwise
# drop meaningful trailing zeros or switch to scientific notation
# (e.g. "5.300" -> 5.3, "0.00001098" -> 1.098E-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '22.379.67'
$ws.Range('E2').Value = '  -0.37%  '
$ws.Range('D3').Value = '1.559.91'
$ws.Range('E3').Value = '  -0.84%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''1.002'
$ws.Range('E5').Value = '  +0.01%  '
$ws.Range('D6').Value = '''285.68'
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('D7').Value = '''0.3645'
$ws.Range('E7').Value = '  -2.77%  '
$ws.Range('D8').Value = '''48.68'
$ws.Range('E8').Value = '  -2.53%  '
$ws.Range('D9').Value = '''0.3332'
$ws.Range('E9').Value = '  -2.21%  '
$ws.Range('D10').Value = '''1.122'
$ws.Range('D11').Value = '''0.07371'
$ws.Range('E11').Value = '  -2.84%  '
$ws.Range('D12').Value = '''1.002'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '''20.69'
$ws.Range('E13').Value = '  -3.33%  '
$ws.Range('D14').Value = '''5.896'
$ws.Range('E14').Value = '  -1.65%  '
$ws.Range('D15').Value = '''6.841'
$ws.Range('E15').Value = '  -1.75%  '
$ws.Range('D16').Value = '1.560.42'
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('D17').Value = '''0.00001098'
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('D18').Value = '''88.59'
$ws.Range('E18').Value = '  -2.87%  '
$ws.Range('D19').Value = '''0.06722'
$ws.Range('E19').Value = '  -0.22%  '
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '''6.290'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').Value = '''15.96'
$ws.Range('E22').Value = '  -2.96%  '
$ws.Range('D23').Value = '''11.89'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').Value = '22.375.57'
$ws.Range('E24').Value = '  -0.35%  '
$ws.Range('D25').Value = '''2.395'
$ws.Range('E25').Value = '  +2.88%  '
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').Value = '''149.17'
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').Value = '''19.30'
$ws.Range('E28').Value = '  -4.33%  '
$ws.Range('D29').Value = '''5.013'
$ws.Range('E29').Value = '  +0.16%  '
$ws.Range('D30').Value = '''122.64'
$ws.Range('E30').Value = '  -2.75%  '
$ws.Range('D31').Value = '1.737.80'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('D32').Value = '''1.054'
$ws.Range('E32').Value = '  +1.05%  '
$ws.Range('D33').Value = '''6.084'
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('D35').Value = '''9.545'
$ws.Range('E35').Value = '  -3.58%  '
$ws.Range('D36').Value = '''0.08217'
$ws.Range('E36').Value = '  -2.81%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02366'
$ws.Range('E37').Value = '  -4.19%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.299'
$ws.Range('E38').Value = '  -6.35%  '
$ws.Range('D39').Value = '''0.2206'
$ws.Range('E39').Value = '  -4.02%  '
$ws.Range('D40').Value = '''0.06342'
$ws.Range('E40').Value = '  -3.25%  '
$ws.Range('D41').Value = '''5.300'
$ws.Range('E41').Value = '  -3.56%  '
$ws.Range('E42').Value = '  -3.07%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = '''0.6020'
$ws.Range('E44').Value = '  -4.43%  '
$ws.Range('D45').Value = '''13.58'
$ws.Range('E45').Value = '  -3.19%  '
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('D47').Value = '''0.5699'
$ws.Range('E47').Value = '  -3.17%  '
$ws.Range('D48').Value = '''124.38'
$ws.Range('E48').Value = '  -4.46%  '
$ws.Range('D49').Value = '''1.998'
$ws.Range('E49').Value = '  -4.72%  '
$ws.Range('E50').Value = '  -2.34%  '
$ws.Range('D51').Value = '''0.07211'
$ws.Range('E51').Value = '  -1.74%  '
